# New weekly price record for "Femacal de La Calera - Ciboulette".
# A new row is inserted right after the header block (at row 19), pushing the
# existing historical rows down by one. The new row re-uses the values of the
# most recent (now last, after shifting) row in the table, and only the date
# (column D) is updated to the new week's date - matching how this weekly
# series is normally appended to.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 19; rows 19-218 shift down to 20-219.
$ws.Rows("19").Insert()

# After the shift, the row that used to be the last data row (218) is now at
# row 219. Copy its values into the freshly inserted row 19.
$lastCol = 18
for ($col = 1; $col -le $lastCol; $col++) {
    $srcCell = $ws.Cells.Item(219, $col)
    $dstCell = $ws.Cells.Item(19, $col)
    $dstCell.Value = $srcCell.Value()
}

# New row gets this week's date (2021-11-30), one week after the previous
# newest date (2021-11-23 / serial 44523).
$ws.Range("D19").Value = 44530
